# envois des mails, générations de fichiers PDF
#
# This script reproduces the target commit:
#  - sheet "Test CU 2 Visiter le site": rows 2..6 no longer record the
#    "Fonctionnalité nécessitant un déploiement" result in column C; their
#    column B result is flipped from "Ko" (red) to "Ok" (green), same as
#    row 7. That makes the shared string "Fonctionnalité nécessitant un
#    déploiement" completely unused, so it disappears from the workbook
#    and every later shared string shifts down by one slot.
#  - selections/active cells are refreshed on several sheets
#  - the active/selected worksheet tab moves from "Faire un don" (sheet 7)
#    back to "Devenir membre" (sheet 1)

$wb = $excel.ActiveWorkbook

$OkColor = 5287936   # RGB(0,176,80)  -> same font color used for "Ok" cells

# --- Sheet 2: "Test CU 2 Visiter le site" ---------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B2").Value = "Ok"
$ws2.Range("B2").Font.Color = $OkColor
$ws2.Range("C2").Value = $null

$ws2.Range("B3").Value = "Ok"
$ws2.Range("B3").Font.Color = $OkColor
$ws2.Range("C3").Value = $null

$ws2.Range("B4").Value = "Ok"
$ws2.Range("B4").Font.Color = $OkColor
$ws2.Range("C4").Value = $null

$ws2.Range("B5").Value = "Ok"
$ws2.Range("B5").Font.Color = $OkColor
$ws2.Range("C5").Value = $null

$ws2.Range("B6").Value = "Ok"
$ws2.Range("B6").Font.Color = $OkColor
$ws2.Range("C6").Value = $null

$ws2.Range("A18").Select() | Out-Null

# --- Sheet 3: "Test CU 3 s'authentifier" -----------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B30").Select() | Out-Null

# --- Sheet 6: "Test CU 6 S'inscrire à la news" ------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("D29").Select() | Out-Null

# --- Sheet 7: "Test CU 7 Faire un don" --------------------------------------
# Its tab was selected before; it is no longer the active tab afterwards,
# but its own remembered selection is unchanged.

# --- Sheet 1: "Test CU 1 Devenir membre" ------------------------------------
# Becomes the active tab again, with a new remembered selection.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A21").Select() | Out-Null
